$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp.
$ws.Name = "IClientBalance-20240814-104249-"

# Column G ("Dt. Referencia") holds the reference date as a serial number.
# Every data row (2-274) moves from 2024-08-13 (45517) to 2024-08-14 (45518).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45518
}

# A handful of rows also got updated balances (Vl. Projetado / Saldo Previsto / Vl. Total).
$ws.Cells.Item(15, 5).Value = 2353.73
$ws.Cells.Item(15, 8).Value = 2353.73

$ws.Cells.Item(43, 5).Value = 1660.96
$ws.Cells.Item(43, 8).Value = 1660.96

$ws.Cells.Item(52, 5).Value = 2104.56
$ws.Cells.Item(52, 8).Value = 2104.56

$ws.Cells.Item(57, 5).Value = 2260.11
$ws.Cells.Item(57, 8).Value = 2260.11

$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 14116.42
$ws.Cells.Item(58, 8).Value = 14116.42

$ws.Cells.Item(101, 5).Value = 1147.49
$ws.Cells.Item(101, 8).Value = 1147.49

$ws.Cells.Item(104, 5).Value = 1676.01
$ws.Cells.Item(104, 8).Value = 1676.01

$ws.Cells.Item(112, 5).Value = 1461.32
$ws.Cells.Item(112, 8).Value = 1461.32

$ws.Cells.Item(113, 5).Value = 1538.45
$ws.Cells.Item(113, 8).Value = 1538.45

$ws.Cells.Item(118, 5).Value = 303.44
$ws.Cells.Item(118, 8).Value = 303.44

$ws.Cells.Item(138, 5).Value = 1217.9000000000001
$ws.Cells.Item(138, 8).Value = 1217.9000000000001

$ws.Cells.Item(143, 5).Value = 2668.91
$ws.Cells.Item(143, 8).Value = 2668.91

$ws.Cells.Item(165, 5).Value = 2015.09
$ws.Cells.Item(165, 8).Value = 2015.09

$ws.Cells.Item(232, 5).Value = 1370.9
$ws.Cells.Item(232, 8).Value = 1370.9

$ws.Cells.Item(255, 5).Value = 1625.51
$ws.Cells.Item(255, 8).Value = 1625.51

$ws.Cells.Item(270, 4).Value = 0
$ws.Cells.Item(270, 5).Value = 943.18

# Reflect the author's final active selection on the sheet.
[void]$ws.Range("B3").Select()
